# Backlog_4.xlsx edit:
# The "Semana" column (C) values were changed from the text label "Semana 04"
# to the plain number 4 on both the SPN and ITI sheets. Once every usage of
# the "Semana 04" shared string is gone, saving drops it from the shared
# string table, which is what shifts all the other <v> string indices seen
# in the diff - that happens automatically, so only the actual value edits
# are needed here.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("SPN")
for ($r = 2; $r -le 30; $r++) {
    $ws1.Cells.Item($r, 3).Value = 4
}

$ws2 = $wb.Worksheets.Item("ITI")
for ($r = 2; $r -le 45; $r++) {
    $ws2.Cells.Item($r, 3).Value = 4
}

# View state: ITI was reviewed (scrolled to its bottom rows, column C
# selected) but ends up not being the active tab when the file is saved.
$ws2.Activate()
[void]$ws2.Range("C2:C45").Select()

# SPN is the last sheet touched, so it ends up the active tab, with its own
# column C selection.
$ws1.Activate()
[void]$ws1.Range("C2:C30").Select()
